$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 30704050
$ws.Range("J70").Value = 25643382
$ws.Range("L70").Value = 76930146
$ws.Range("N70").Value = -76930686
$ws.Range("H73").Value = 30704050
$ws.Range("J73").Value = 25643382
$ws.Range("L73").Value = 76930146
$ws.Range("N73").Value = -76932018
$ws.Range("H103").Value = 932.2632
$ws.Range("J103").Value = 1012.1875
$ws.Range("L103").Value = 3036.5625
$ws.Range("N103").Value = -4208.5625
$ws.Range("H125").Value = 62500856
$ws.Range("I125").Value = 250000350
$ws.Range("J125").Value = 1025
$ws.Range("K125").Value = 2250003150
$ws.Range("L125").Value = 9225
$ws.Range("M125").Value = -2250000690
$ws.Range("N125").Value = -14145
$ws.Range("H137").Value = 3835.5833
$ws.Range("I137").Value = 4197.875
$ws.Range("J137").Value = 3111
$ws.Range("K137").Value = 12593.625
$ws.Range("L137").Value = 9333
$ws.Range("M137").Value = -10043.625
$ws.Range("N137").Value = -14433

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1988312.2
$ws.Range("I32").Value = 2020371.2
$ws.Range("K32").Value = 2020371.2
$ws.Range("M32").Value = -2020084.2
$ws.Range("H45").Value = 5834
$ws.Range("I45").Value = 2861.8572
$ws.Range("K45").Value = 2861.8572
$ws.Range("M45").Value = -2484.8572
$ws.Range("H74").Value = 41313.223
$ws.Range("I74").Value = 74390.28999999999
$ws.Range("J74").Value = 5691.769
$ws.Range("K74").Value = 74390.28999999999
$ws.Range("L74").Value = 5691.769
$ws.Range("M74").Value = -73516.28999999999
$ws.Range("N74").Value = -7439.769
$ws.Range("H77").Value = 41313.223
$ws.Range("I77").Value = 74390.28999999999
$ws.Range("J77").Value = 5691.769
$ws.Range("K77").Value = 371951.45
$ws.Range("L77").Value = 28458.845
$ws.Range("M77").Value = -367583.45
$ws.Range("N77").Value = -37194.845
$ws.Range("H102").Value = 1501
$ws.Range("I102").Value = 1501
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1501
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 121
$ws.Range("H122").Value = 6070.5713
$ws.Range("I122").Value = 3246.5
$ws.Range("J122").Value = 7200.2
$ws.Range("K122").Value = 9739.5
$ws.Range("L122").Value = 21600.6
$ws.Range("M122").Value = -7289.5
$ws.Range("N122").Value = -26500.6
$ws.Range("N102").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("H134").Value = 6522.231
$ws.Range("I134").Value = 2593.9412
$ws.Range("K134").Value = 7781.823600000001
$ws.Range("M134").Value = -5246.823600000001
$ws.Range("N9").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13344540
$ws.Range("I31").Value = 6138.6
$ws.Range("J31").Value = 16679140
$ws.Range("K31").Value = 6138.6
$ws.Range("L31").Value = 16679140
$ws.Range("M31").Value = -5843.6
$ws.Range("N31").Value = -16679730
$ws.Range("H34").Value = 13344540
$ws.Range("I34").Value = 6138.6
$ws.Range("J34").Value = 16679140
$ws.Range("K34").Value = 6138.6
$ws.Range("L34").Value = 16679140
$ws.Range("M34").Value = -5936.6
$ws.Range("N34").Value = -16679544
$ws.Range("H94").Value = 932.05884
$ws.Range("I94").Value = 1237.375
$ws.Range("J94").Value = 660.6667
$ws.Range("K94").Value = 1237.375
$ws.Range("L94").Value = 660.6667
$ws.Range("M94").Value = -786.375
$ws.Range("N94").Value = -1562.6667
$ws.Range("H107").Value = 2450.2727
$ws.Range("I107").Value = 1813.6923
$ws.Range("K107").Value = 1813.6923
$ws.Range("M107").Value = 106.3077000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 288385340
$ws.Range("I4").Value = 336667330
$ws.Range("J4").Value = 252173820
$ws.Range("K4").Value = 1010001990
$ws.Range("L4").Value = 756521460
$ws.Range("M4").Value = -1010001878
$ws.Range("N4").Value = -756521684
$ws.Range("H34").Value = 5923.4614
$ws.Range("J34").Value = 6996.636
$ws.Range("L34").Value = 20989.908
$ws.Range("N34").Value = -21157.908
$ws.Range("H39").Value = 8777.571
$ws.Range("J39").Value = 12054.556
$ws.Range("L39").Value = 36163.66800000001
$ws.Range("N39").Value = -36751.66800000001
$ws.Range("H55").Value = 75443800
$ws.Range("J55").Value = 9099712
$ws.Range("L55").Value = 27299136
$ws.Range("N55").Value = -27299490
$ws.Range("H60").Value = 728.1667
$ws.Range("I60").Value = 184.5
$ws.Range("J60").Value = 1000
$ws.Range("K60").Value = 553.5
$ws.Range("L60").Value = 3000
$ws.Range("M60").Value = -302.5
$ws.Range("N60").Value = -3502
$ws.Range("H64").Value = 100000000
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("H67").Value = 100000000
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("H87").Value = 500000500
$ws.Range("I87").Value = 500000500
$ws.Range("K87").Value = 1500001500
$ws.Range("M87").Value = -1500000252
$ws.Range("H90").Value = 500000500
$ws.Range("I90").Value = 500000500
$ws.Range("K90").Value = 4500004500
$ws.Range("M90").Value = -4499998260
$ws.Range("H137").Value = 101964.75
$ws.Range("I137").Value = 84777.336
$ws.Range("J137").Value = 127745.875
$ws.Range("K137").Value = 254332.008
$ws.Range("L137").Value = 383237.625
$ws.Range("M137").Value = -249232.008
$ws.Range("N137").Value = -393437.625
$ws.Range("M64").ClearContents()
$ws.Range("M67").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2685960.2
$ws.Range("I122").Value = 4026553.2
$ws.Range("J122").Value = 4774.1113
$ws.Range("K122").Value = 12079659.6
$ws.Range("L122").Value = 14322.3339
$ws.Range("M122").Value = -12077209.6
$ws.Range("N122").Value = -19222.3339
$ws.Range("H132").Value = 6760.9165
$ws.Range("I132").Value = 1324.75
$ws.Range("J132").Value = 9479
$ws.Range("K132").Value = 3974.25
$ws.Range("L132").Value = 28437
$ws.Range("M132").Value = -1444.25
$ws.Range("N132").Value = -33497

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5500.68
$ws.Range("I7").Value = 4035.7058
$ws.Range("J7").Value = 8613.75
$ws.Range("K7").Value = 4035.7058
$ws.Range("L7").Value = 8613.75
$ws.Range("M7").Value = -3923.7058
$ws.Range("N7").Value = -8837.75
$ws.Range("H22").Value = 2023.4445
$ws.Range("I22").Value = 900
$ws.Range("K22").Value = 900
$ws.Range("M22").Value = -605
$ws.Range("H27").Value = 2023.4445
$ws.Range("I27").Value = 900
$ws.Range("K27").Value = 900
$ws.Range("M27").Value = -793
$ws.Range("H40").Value = 4339.375
$ws.Range("I40").Value = 3835
$ws.Range("J40").Value = 4684.4736
$ws.Range("K40").Value = 3835
$ws.Range("L40").Value = 4684.4736
$ws.Range("M40").Value = -3699
$ws.Range("N40").Value = -4956.4736
$ws.Range("H46").Value = 1500803.9
$ws.Range("J46").Value = 2401.6155
$ws.Range("L46").Value = 2401.6155
$ws.Range("N46").Value = -2777.6155
$ws.Range("H87").Value = 70000
$ws.Range("J87").Value = 70000
$ws.Range("L87").Value = 70000
$ws.Range("N87").Value = -72246
$ws.Range("H90").Value = 70000
$ws.Range("J90").Value = 70000
$ws.Range("L90").Value = 210000
$ws.Range("N90").Value = -221232
$ws.Range("H93").Value = 5322.0625
$ws.Range("J93").Value = 10381.4
$ws.Range("L93").Value = 10381.4
$ws.Range("N93").Value = -12877.4
$ws.Range("H122").Value = 3425.0356
$ws.Range("I122").Value = 2602.8096
$ws.Range("J122").Value = 5891.7144
$ws.Range("K122").Value = 7808.4288
$ws.Range("L122").Value = 17675.1432
$ws.Range("M122").Value = -5358.4288
$ws.Range("N122").Value = -22575.1432
$ws.Range("H126").Value = 5500.68
$ws.Range("I126").Value = 4035.7058
$ws.Range("J126").Value = 8613.75
$ws.Range("K126").Value = 12107.1174
$ws.Range("L126").Value = 25841.25
$ws.Range("M126").Value = -9637.117400000001
$ws.Range("N126").Value = -30781.25
$ws.Range("H132").Value = 11911956
$ws.Range("I132").Value = 26319114
$ws.Range("J132").Value = 10391.261
$ws.Range("K132").Value = 78957342
$ws.Range("L132").Value = 31173.783
$ws.Range("M132").Value = -78954812
$ws.Range("N132").Value = -36233.783
$ws.Range("H136").Value = 11259.56
$ws.Range("I136").Value = 3697.8
$ws.Range("K136").Value = 11093.4
$ws.Range("M136").Value = -8543.400000000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2359.5293
$ws.Range("J113").Value = 3435.3333
$ws.Range("L113").Value = 10305.9999
$ws.Range("N113").Value = -14645.9999
$ws.Range("H122").Value = 337930.9
$ws.Range("I122").Value = 669235.2
$ws.Range("J122").Value = 6626.6665
$ws.Range("K122").Value = 2007705.6
$ws.Range("L122").Value = 19879.9995
$ws.Range("M122").Value = -2005255.6
$ws.Range("N122").Value = -24779.9995
$ws.Range("H132").Value = 41670616
$ws.Range("I132").Value = 250001500
$ws.Range("J132").Value = 4438.7
$ws.Range("K132").Value = 750004500
$ws.Range("L132").Value = 13316.1
$ws.Range("M132").Value = -750001970
$ws.Range("N132").Value = -18376.1
$ws.Range("H136").Value = 35755930
$ws.Range("J136").Value = 54996.57
$ws.Range("L136").Value = 164989.71
$ws.Range("N136").Value = -170089.71
